# Feature : 근접 무기 Sword 추가
# Insert a new "Sword" melee weapon row right after the Axe rows (becomes
# row 5), pushing the existing Spear/Bow/Wand rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at row 5 - everything that used to be row 5 (Spear),
#    6 (Bow), 7 (Wand) shifts down to 6, 7, 8.
$ws.Rows(5).Insert()

# 2) Populate the new row 5 with the Sword entry.
$ws.Range("A5").Value = 10112011
$ws.Range("B5").Value = "Weapon"
$ws.Range("C5").Value = "Sword"
$ws.Range("D5").Value = "칼"
$ws.Range("E5").Value = "Items/Icons/Weapons/Melee/Sword_1"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = "Items/Prefabs/Weapons/Melee/Sword_1"
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = "Normal"
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Melee"
$ws.Range("S5").Value = "None"
$ws.Range("T5").Value = 0

# 3) The row that used to be row 5 (Spear) is now row 6 - fix up its id.
$ws.Range("A6").Value = 10113011

# 4) Column width tweaks that came with the edit.
$ws.Columns(1).ColumnWidth = 12.857
$ws.Columns(4).ColumnWidth = 19.571

# 5) Selection left on M5 by the editing author.
$null = $ws.Range("M5").Select()
